# Restructured Request, Changed Student
# Adds a new "Assigned" boolean column (D) to the student list sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column D
$ws.Range("D1").Value = "Assigned"

# Row 2 student is now "Assigned" (TRUE); the rest default to "not assigned" (FALSE)
$ws.Range("D2").Value = $true
$ws.Range("D3").Value = $false
$ws.Range("D4").Value = $false
$ws.Range("D5").Value = $false
$ws.Range("D6").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("D8").Value = $false
$ws.Range("D9").Value = $false
$ws.Range("D10").Value = $false
$ws.Range("D11").Value = $false
$ws.Range("D12").Value = $false

# Re-apply the default font so Excel records an explicit (but visually identical)
# font style for the new boolean column, matching fontId 0 with applyFont set.
$ws.Range("D2:D12").Font.Bold = $false

# Update the active selection to match the authored change
$ws.Range("D3").Select()
